# params.xlsx: update parameter rows and add the 21.03.2020 measurement row.
#
# Resulting "date" column:
#   row2 27.02.2020 (unchanged date, D changes)
#   row3 17.03.2020 -> 07.03.2020 (date corrected, parameters updated)
#   row4 21.03.2020 -> 10.03.2020 (new measurement date, parameters updated)
#   row5 (new)      -> 21.03.2020 (re-added as a new row with its own parameters)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date-like labels in column A must stay as plain text (shared strings),
# not get auto-converted into Excel date serials. Temporarily force a text
# number format while assigning them, then restore "General" so the cell
# style stays the workbook's default.
$dateCol = $ws.Range("A3:A5")
$dateCol.NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "07.03.2020"
$ws.Cells.Item(4,1).Value = "10.03.2020"
$ws.Cells.Item(5,1).Value = "21.03.2020"
$dateCol.NumberFormat = "General"

# Row 2 (27.02.2020) - only the mpic parameter (column D) changed.
$ws.Cells.Item(2,4).Value = 0.2

# Row 3 (07.03.2020) - mlam/mpic/mlag/vlag/vlos updated.
$ws.Cells.Item(3,2).Value = 1.25
$ws.Cells.Item(3,3).Value = 0.01
$ws.Cells.Item(3,4).Value = 0.15
$ws.Cells.Item(3,5).Value = 0.05
$ws.Cells.Item(3,6).Value = 8
$ws.Cells.Item(3,7).Value = 12
$ws.Cells.Item(3,8).Value = 9
$ws.Cells.Item(3,9).Value = 15

# Row 4 (10.03.2020) - same new parameter set as row 3 except mpic.
$ws.Cells.Item(4,2).Value = 1.25
$ws.Cells.Item(4,3).Value = 0.01
$ws.Cells.Item(4,4).Value = 0.05
$ws.Cells.Item(4,5).Value = 0.05
$ws.Cells.Item(4,6).Value = 8
$ws.Cells.Item(4,7).Value = 12
$ws.Cells.Item(4,8).Value = 9
$ws.Cells.Item(4,9).Value = 15

# Row 5 (21.03.2020, newly appended) - reuses the old mlam/vlam baseline values.
$ws.Cells.Item(5,2).Value = 1.12
$ws.Cells.Item(5,3).Value = 0.01
$ws.Cells.Item(5,4).Value = 0.05
$ws.Cells.Item(5,5).Value = 0.05
$ws.Cells.Item(5,6).Value = 8
$ws.Cells.Item(5,7).Value = 12
$ws.Cells.Item(5,8).Value = 9
$ws.Cells.Item(5,9).Value = 15

# Selection moves to D3, matching the saved view state.
$ws.Range("D3").Select()
